$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Agrn"
$ws.Cells.Item(2, 3).Value = "Atp1a3"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 17.05375133333333
$ws.Cells.Item(2, 8).Value = 51.161254
$ws.Cells.Item(2, 9).Value = 0.3501286198398134
$ws.Cells.Item(2, 10).Value = 0.3501286198398134
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.012411
$ws.Cells.Item(2, 14).Value = 0.037233
$ws.Cells.Item(2, 15).Value = 0.0001734784643863198
$ws.Cells.Item(2, 16).Value = 0.0001734784643863198
$ws.Cells.Item(2, 17).Value = 0.211654107798
$ws.Cells.Item(2, 18).Value = 1.904886970182
$ws.Cells.Item(2, 19).Value = 0.00006073977530751236
$ws.Cells.Item(2, 20).Value = 0.00006073977530751239

$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Agrn"
$ws.Cells.Item(3, 3).Value = "Atp1a3"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 17.05375133333333
$ws.Cells.Item(3, 8).Value = 51.161254
$ws.Cells.Item(3, 9).Value = 0.3501286198398134
$ws.Cells.Item(3, 10).Value = 0.3501286198398134
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 0.08803833333333333
$ws.Cells.Item(3, 14).Value = 0.264115
$ws.Cells.Item(3, 15).Value = 0.001230582134702894
$ws.Cells.Item(3, 16).Value = 0.001230582134702894
$ws.Cells.Item(3, 17).Value = 1.501383844467778
$ws.Cells.Item(3, 18).Value = 13.51245460021
$ws.Cells.Item(3, 19).Value = 0.0004308620244230555
$ws.Cells.Item(3, 20).Value = 0.0004308620244230557

$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Agrn"
$ws.Cells.Item(4, 3).Value = "Atp1a3"
$ws.Cells.Item(4, 4).Value = "M1"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 17.05375133333333
$ws.Cells.Item(4, 8).Value = 51.161254
$ws.Cells.Item(4, 9).Value = 0.3501286198398134
$ws.Cells.Item(4, 10).Value = 0.3501286198398134
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 20.739774
$ws.Cells.Item(4, 14).Value = 62.219322
$ws.Cells.Item(4, 15).Value = 0.2898963939440272
$ws.Cells.Item(4, 16).Value = 0.2898963939440272
$ws.Cells.Item(4, 17).Value = 353.690948505532
$ws.Cells.Item(4, 18).Value = 3183.218536549788
$ws.Cells.Item(4, 19).Value = 0.1015010243081611
$ws.Cells.Item(4, 20).Value = 0.1015010243081611

$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Agrn"
$ws.Cells.Item(5, 3).Value = "Atp1a3"
$ws.Cells.Item(5, 4).Value = "M2"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 17.05375133333333
$ws.Cells.Item(5, 8).Value = 51.161254
$ws.Cells.Item(5, 9).Value = 0.3501286198398134
$ws.Cells.Item(5, 10).Value = 0.3501286198398134
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 19.56635666666667
$ws.Cells.Item(5, 14).Value = 58.69907
$ws.Cells.Item(5, 15).Value = 0.2734946022212847
$ws.Cells.Item(5, 16).Value = 0.2734946022212847
$ws.Cells.Item(5, 17).Value = 333.6797810926423
$ws.Cells.Item(5, 18).Value = 3003.11802983378
$ws.Cells.Item(5, 19).Value = 0.09575828760937716
$ws.Cells.Item(5, 20).Value = 0.09575828760937717

$ws.Cells.Item(6, 1).Value = "ECs"
$ws.Cells.Item(6, 2).Value = "Agrn"
$ws.Cells.Item(6, 3).Value = "Atp1a3"
$ws.Cells.Item(6, 4).Value = "Neutro"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 17.05375133333333
$ws.Cells.Item(6, 8).Value = 51.161254
$ws.Cells.Item(6, 9).Value = 0.3501286198398134
$ws.Cells.Item(6, 10).Value = 0.3501286198398134
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 31.06110666666667
$ws.Cells.Item(6, 14).Value = 93.18332000000001
$ws.Cells.Item(6, 15).Value = 0.4341659082002267
$ws.Cells.Item(6, 16).Value = 0.4341659082002268
$ws.Cells.Item(6, 17).Value = 529.7083892314756
$ws.Cells.Item(6, 18).Value = 4767.375503083281
$ws.Cells.Item(6, 19).Value = 0.1520139102196445
$ws.Cells.Item(6, 20).Value = 0.1520139102196446

$ws.Cells.Item(7, 1).Value = "ECs"
$ws.Cells.Item(7, 2).Value = "Agrn"
$ws.Cells.Item(7, 3).Value = "Atp1a3"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 17.05375133333333
$ws.Cells.Item(7, 8).Value = 51.161254
$ws.Cells.Item(7, 9).Value = 0.3501286198398134
$ws.Cells.Item(7, 10).Value = 0.3501286198398134
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.07433466666666666
$ws.Cells.Item(7, 14).Value = 0.223004
$ws.Cells.Item(7, 15).Value = 0.001039035035372032
$ws.Cells.Item(7, 16).Value = 0.001039035035372032
$ws.Cells.Item(7, 17).Value = 1.267684920779556
$ws.Cells.Item(7, 18).Value = 11.409164287016
$ws.Cells.Item(7, 19).Value = 0.0003637959029000211
$ws.Cells.Item(7, 20).Value = 0.0003637959029000212

$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Agrn"
$ws.Cells.Item(8, 3).Value = "Atp1a3"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 3.037017666666667
$ws.Cells.Item(8, 8).Value = 9.111053
$ws.Cells.Item(8, 9).Value = 0.06235266266494154
$ws.Cells.Item(8, 10).Value = 0.06235266266494155
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 0.012411
$ws.Cells.Item(8, 14).Value = 0.037233
$ws.Cells.Item(8, 15).Value = 0.0001734784643863198
$ws.Cells.Item(8, 16).Value = 0.0001734784643863198
$ws.Cells.Item(8, 17).Value = 0.037692426261
$ws.Cells.Item(8, 18).Value = 0.339231836349
$ws.Cells.Item(8, 19).Value = 0.00001081684416951227
$ws.Cells.Item(8, 20).Value = 0.00001081684416951228

$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Agrn"
$ws.Cells.Item(9, 3).Value = "Atp1a3"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 3.037017666666667
$ws.Cells.Item(9, 8).Value = 9.111053
$ws.Cells.Item(9, 9).Value = 0.06235266266494154
$ws.Cells.Item(9, 10).Value = 0.06235266266494155
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.08803833333333333
$ws.Cells.Item(9, 14).Value = 0.264115
$ws.Cells.Item(9, 15).Value = 0.001230582134702894
$ws.Cells.Item(9, 16).Value = 0.001230582134702894
$ws.Cells.Item(9, 17).Value = 0.2673739736772222
$ws.Cells.Item(9, 18).Value = 2.406365763095
$ws.Cells.Item(9, 19).Value = 0.00007673007272663318
$ws.Cells.Item(9, 20).Value = 0.00007673007272663322

$ws.Cells.Item(10, 1).Value = "FAPs"
$ws.Cells.Item(10, 2).Value = "Agrn"
$ws.Cells.Item(10, 3).Value = "Atp1a3"
$ws.Cells.Item(10, 4).Value = "M1"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 3.037017666666667
$ws.Cells.Item(10, 8).Value = 9.111053
$ws.Cells.Item(10, 9).Value = 0.06235266266494154
$ws.Cells.Item(10, 10).Value = 0.06235266266494155
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 20.739774
$ws.Cells.Item(10, 14).Value = 62.219322
$ws.Cells.Item(10, 15).Value = 0.2898963939440272
$ws.Cells.Item(10, 16).Value = 0.2898963939440272
$ws.Cells.Item(10, 17).Value = 62.987060040674
$ws.Cells.Item(10, 18).Value = 566.883540366066
$ws.Cells.Item(10, 19).Value = 0.01807581205937493
$ws.Cells.Item(10, 20).Value = 0.01807581205937493

$ws.Cells.Item(11, 1).Value = "FAPs"
$ws.Cells.Item(11, 2).Value = "Agrn"
$ws.Cells.Item(11, 3).Value = "Atp1a3"
$ws.Cells.Item(11, 4).Value = "M2"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 3.037017666666667
$ws.Cells.Item(11, 8).Value = 9.111053
$ws.Cells.Item(11, 9).Value = 0.06235266266494154
$ws.Cells.Item(11, 10).Value = 0.06235266266494155
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 19.56635666666667
$ws.Cells.Item(11, 14).Value = 58.69907
$ws.Cells.Item(11, 15).Value = 0.2734946022212847
$ws.Cells.Item(11, 16).Value = 0.2734946022212847
$ws.Cells.Item(11, 17).Value = 59.42337086896779
$ws.Cells.Item(11, 18).Value = 534.8103378207099
$ws.Cells.Item(11, 19).Value = 0.01705311667298614
$ws.Cells.Item(11, 20).Value = 0.01705311667298614

$ws.Cells.Item(12, 1).Value = "FAPs"
$ws.Cells.Item(12, 2).Value = "Agrn"
$ws.Cells.Item(12, 3).Value = "Atp1a3"
$ws.Cells.Item(12, 4).Value = "Neutro"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 3.037017666666667
$ws.Cells.Item(12, 8).Value = 9.111053
$ws.Cells.Item(12, 9).Value = 0.06235266266494154
$ws.Cells.Item(12, 10).Value = 0.06235266266494155
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 31.06110666666667
$ws.Cells.Item(12, 14).Value = 93.18332000000001
$ws.Cells.Item(12, 15).Value = 0.4341659082002267
$ws.Cells.Item(12, 16).Value = 0.4341659082002268
$ws.Cells.Item(12, 17).Value = 94.33312969288447
$ws.Cells.Item(12, 18).Value = 848.9981672359601
$ws.Cells.Item(12, 19).Value = 0.02707140041462671
$ws.Cells.Item(12, 20).Value = 0.02707140041462672

$ws.Cells.Item(13, 1).Value = "FAPs"
$ws.Cells.Item(13, 2).Value = "Agrn"
$ws.Cells.Item(13, 3).Value = "Atp1a3"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 3.037017666666667
$ws.Cells.Item(13, 8).Value = 9.111053
$ws.Cells.Item(13, 9).Value = 0.06235266266494154
$ws.Cells.Item(13, 10).Value = 0.06235266266494155
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 0.07433466666666666
$ws.Cells.Item(13, 14).Value = 0.223004
$ws.Cells.Item(13, 15).Value = 0.001039035035372032
$ws.Cells.Item(13, 16).Value = 0.001039035035372032
$ws.Cells.Item(13, 17).Value = 0.2257556959124444
$ws.Cells.Item(13, 18).Value = 2.031801263212
$ws.Cells.Item(13, 19).Value = 0.00006478660105760789
$ws.Cells.Item(13, 20).Value = 0.00006478660105760791

$ws.Cells.Item(14, 1).Value = "M1"
$ws.Cells.Item(14, 2).Value = "Agrn"
$ws.Cells.Item(14, 3).Value = "Atp1a3"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 4.955296333333334
$ws.Cells.Item(14, 8).Value = 14.865889
$ws.Cells.Item(14, 9).Value = 0.101736622762645
$ws.Cells.Item(14, 10).Value = 0.101736622762645
$ws.Cells.Item(14, 11).Value = 2
$ws.Cells.Item(14, 12).Value = 0.6666666666666666
$ws.Cells.Item(14, 13).Value = 0.012411
$ws.Cells.Item(14, 14).Value = 0.037233
$ws.Cells.Item(14, 15).Value = 0.0001734784643863198
$ws.Cells.Item(14, 16).Value = 0.0001734784643863198
$ws.Cells.Item(14, 17).Value = 0.06150018279300001
$ws.Cells.Item(14, 18).Value = 0.5535016451370001
$ws.Cells.Item(14, 19).Value = 0.00001764911308871396
$ws.Cells.Item(14, 20).Value = 0.00001764911308871397

$ws.Cells.Item(15, 1).Value = "M1"
$ws.Cells.Item(15, 2).Value = "Agrn"
$ws.Cells.Item(15, 3).Value = "Atp1a3"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 4.955296333333334
$ws.Cells.Item(15, 8).Value = 14.865889
$ws.Cells.Item(15, 9).Value = 0.101736622762645
$ws.Cells.Item(15, 10).Value = 0.101736622762645
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 0.08803833333333333
$ws.Cells.Item(15, 14).Value = 0.264115
$ws.Cells.Item(15, 15).Value = 0.001230582134702894
$ws.Cells.Item(15, 16).Value = 0.001230582134702894
$ws.Cells.Item(15, 17).Value = 0.4362560303594445
$ws.Cells.Item(15, 18).Value = 3.926304273235
$ws.Cells.Item(15, 19).Value = 0.0001251952704167187
$ws.Cells.Item(15, 20).Value = 0.0001251952704167188

$ws.Cells.Item(16, 1).Value = "M1"
$ws.Cells.Item(16, 2).Value = "Agrn"
$ws.Cells.Item(16, 3).Value = "Atp1a3"
$ws.Cells.Item(16, 4).Value = "M1"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 4.955296333333334
$ws.Cells.Item(16, 8).Value = 14.865889
$ws.Cells.Item(16, 9).Value = 0.101736622762645
$ws.Cells.Item(16, 10).Value = 0.101736622762645
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 20.739774
$ws.Cells.Item(16, 14).Value = 62.219322
$ws.Cells.Item(16, 15).Value = 0.2898963939440272
$ws.Cells.Item(16, 16).Value = 0.2898963939440272
$ws.Cells.Item(16, 17).Value = 102.771726056362
$ws.Cells.Item(16, 18).Value = 924.945534507258
$ws.Cells.Item(16, 19).Value = 0.02949308007093462
$ws.Cells.Item(16, 20).Value = 0.02949308007093463

$ws.Cells.Item(17, 1).Value = "M1"
$ws.Cells.Item(17, 2).Value = "Agrn"
$ws.Cells.Item(17, 3).Value = "Atp1a3"
$ws.Cells.Item(17, 4).Value = "M2"
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 4.955296333333334
$ws.Cells.Item(17, 8).Value = 14.865889
$ws.Cells.Item(17, 9).Value = 0.101736622762645
$ws.Cells.Item(17, 10).Value = 0.101736622762645
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 19.56635666666667
$ws.Cells.Item(17, 14).Value = 58.69907
$ws.Cells.Item(17, 15).Value = 0.2734946022212847
$ws.Cells.Item(17, 16).Value = 0.2734946022212847
$ws.Cells.Item(17, 17).Value = 96.95709544702557
$ws.Cells.Item(17, 18).Value = 872.61385902323
$ws.Cells.Item(17, 19).Value = 0.02782441717380649
$ws.Cells.Item(17, 20).Value = 0.0278244171738065

$ws.Cells.Item(18, 1).Value = "M1"
$ws.Cells.Item(18, 2).Value = "Agrn"
$ws.Cells.Item(18, 3).Value = "Atp1a3"
$ws.Cells.Item(18, 4).Value = "Neutro"
$ws.Cells.Item(18, 5).Value = 3
$ws.Cells.Item(18, 6).Value = 1
$ws.Cells.Item(18, 7).Value = 4.955296333333334
$ws.Cells.Item(18, 8).Value = 14.865889
$ws.Cells.Item(18, 9).Value = 0.101736622762645
$ws.Cells.Item(18, 10).Value = 0.101736622762645
$ws.Cells.Item(18, 11).Value = 3
$ws.Cells.Item(18, 12).Value = 1
$ws.Cells.Item(18, 13).Value = 31.06110666666667
$ws.Cells.Item(18, 14).Value = 93.18332000000001
$ws.Cells.Item(18, 15).Value = 0.4341659082002267
$ws.Cells.Item(18, 16).Value = 0.4341659082002268
$ws.Cells.Item(18, 17).Value = 153.9169879746089
$ws.Cells.Item(18, 18).Value = 1385.25289177148
$ws.Cells.Item(18, 19).Value = 0.04417057321896763
$ws.Cells.Item(18, 20).Value = 0.04417057321896765

$ws.Cells.Item(19, 1).Value = "M1"
$ws.Cells.Item(19, 2).Value = "Agrn"
$ws.Cells.Item(19, 3).Value = "Atp1a3"
$ws.Cells.Item(19, 4).Value = "sCs"
$ws.Cells.Item(19, 5).Value = 3
$ws.Cells.Item(19, 6).Value = 1
$ws.Cells.Item(19, 7).Value = 4.955296333333334
$ws.Cells.Item(19, 8).Value = 14.865889
$ws.Cells.Item(19, 9).Value = 0.101736622762645
$ws.Cells.Item(19, 10).Value = 0.101736622762645
$ws.Cells.Item(19, 11).Value = 3
$ws.Cells.Item(19, 12).Value = 1
$ws.Cells.Item(19, 13).Value = 0.07433466666666666
$ws.Cells.Item(19, 14).Value = 0.223004
$ws.Cells.Item(19, 15).Value = 0.001039035035372032
$ws.Cells.Item(19, 16).Value = 0.001039035035372032
$ws.Cells.Item(19, 17).Value = 0.3683503011728889
$ws.Cells.Item(19, 18).Value = 3.315152710556
$ws.Cells.Item(19, 19).Value = 0.0001057079154308159
$ws.Cells.Item(19, 20).Value = 0.0001057079154308159

$ws.Cells.Item(20, 1).Value = "M2"
$ws.Cells.Item(20, 2).Value = "Agrn"
$ws.Cells.Item(20, 3).Value = "Atp1a3"
$ws.Cells.Item(20, 4).Value = "ECs"
$ws.Cells.Item(20, 5).Value = 3
$ws.Cells.Item(20, 6).Value = 1
$ws.Cells.Item(20, 7).Value = 5.897008666666667
$ws.Cells.Item(20, 8).Value = 17.691026
$ws.Cells.Item(20, 9).Value = 0.1210708110659339
$ws.Cells.Item(20, 10).Value = 0.1210708110659339
$ws.Cells.Item(20, 11).Value = 2
$ws.Cells.Item(20, 12).Value = 0.6666666666666666
$ws.Cells.Item(20, 13).Value = 0.012411
$ws.Cells.Item(20, 14).Value = 0.037233
$ws.Cells.Item(20, 15).Value = 0.0001734784643863198
$ws.Cells.Item(20, 16).Value = 0.0001734784643863198
$ws.Cells.Item(20, 17).Value = 0.07318777456200001
$ws.Cells.Item(20, 18).Value = 0.6586899710580001
$ws.Cells.Item(20, 19).Value = 0.00002100317838572446
$ws.Cells.Item(20, 20).Value = 0.00002100317838572447

$ws.Cells.Item(21, 1).Value = "M2"
$ws.Cells.Item(21, 2).Value = "Agrn"
$ws.Cells.Item(21, 3).Value = "Atp1a3"
$ws.Cells.Item(21, 4).Value = "FAPs"
$ws.Cells.Item(21, 5).Value = 3
$ws.Cells.Item(21, 6).Value = 1
$ws.Cells.Item(21, 7).Value = 5.897008666666667
$ws.Cells.Item(21, 8).Value = 17.691026
$ws.Cells.Item(21, 9).Value = 0.1210708110659339
$ws.Cells.Item(21, 10).Value = 0.1210708110659339
$ws.Cells.Item(21, 11).Value = 3
$ws.Cells.Item(21, 12).Value = 1
$ws.Cells.Item(21, 13).Value = 0.08803833333333333
$ws.Cells.Item(21, 14).Value = 0.264115
$ws.Cells.Item(21, 15).Value = 0.001230582134702894
$ws.Cells.Item(21, 16).Value = 0.001230582134702894
$ws.Cells.Item(21, 17).Value = 0.5191628146655556
$ws.Cells.Item(21, 18).Value = 4.67246533199
$ws.Cells.Item(21, 19).Value = 0.0001489875771317277
$ws.Cells.Item(21, 20).Value = 0.0001489875771317277

$ws.Cells.Item(22, 1).Value = "M2"
$ws.Cells.Item(22, 2).Value = "Agrn"
$ws.Cells.Item(22, 3).Value = "Atp1a3"
$ws.Cells.Item(22, 4).Value = "M1"
$ws.Cells.Item(22, 5).Value = 3
$ws.Cells.Item(22, 6).Value = 1
$ws.Cells.Item(22, 7).Value = 5.897008666666667
$ws.Cells.Item(22, 8).Value = 17.691026
$ws.Cells.Item(22, 9).Value = 0.1210708110659339
$ws.Cells.Item(22, 10).Value = 0.1210708110659339
$ws.Cells.Item(22, 11).Value = 3
$ws.Cells.Item(22, 12).Value = 1
$ws.Cells.Item(22, 13).Value = 20.739774
$ws.Cells.Item(22, 14).Value = 62.219322
$ws.Cells.Item(22, 15).Value = 0.2898963939440272
$ws.Cells.Item(22, 16).Value = 0.2898963939440272
$ws.Cells.Item(22, 17).Value = 122.302627022708
$ws.Cells.Item(22, 18).Value = 1100.723643204372
$ws.Cells.Item(22, 19).Value = 0.03509799153989285
$ws.Cells.Item(22, 20).Value = 0.03509799153989286

$ws.Cells.Item(23, 1).Value = "M2"
$ws.Cells.Item(23, 2).Value = "Agrn"
$ws.Cells.Item(23, 3).Value = "Atp1a3"
$ws.Cells.Item(23, 4).Value = "M2"
$ws.Cells.Item(23, 5).Value = 3
$ws.Cells.Item(23, 6).Value = 1
$ws.Cells.Item(23, 7).Value = 5.897008666666667
$ws.Cells.Item(23, 8).Value = 17.691026
$ws.Cells.Item(23, 9).Value = 0.1210708110659339
$ws.Cells.Item(23, 10).Value = 0.1210708110659339
$ws.Cells.Item(23, 11).Value = 3
$ws.Cells.Item(23, 12).Value = 1
$ws.Cells.Item(23, 13).Value = 19.56635666666667
$ws.Cells.Item(23, 14).Value = 58.69907
$ws.Cells.Item(23, 15).Value = 0.2734946022212847
$ws.Cells.Item(23, 16).Value = 0.2734946022212847
$ws.Cells.Item(23, 17).Value = 115.3829748384245
$ws.Cells.Item(23, 18).Value = 1038.44677354582
$ws.Cells.Item(23, 19).Value = 0.0331122133130859
$ws.Cells.Item(23, 20).Value = 0.03311221331308591

$ws.Cells.Item(24, 1).Value = "M2"
$ws.Cells.Item(24, 2).Value = "Agrn"
$ws.Cells.Item(24, 3).Value = "Atp1a3"
$ws.Cells.Item(24, 4).Value = "Neutro"
$ws.Cells.Item(24, 5).Value = 3
$ws.Cells.Item(24, 6).Value = 1
$ws.Cells.Item(24, 7).Value = 5.897008666666667
$ws.Cells.Item(24, 8).Value = 17.691026
$ws.Cells.Item(24, 9).Value = 0.1210708110659339
$ws.Cells.Item(24, 10).Value = 0.1210708110659339
$ws.Cells.Item(24, 11).Value = 3
$ws.Cells.Item(24, 12).Value = 1
$ws.Cells.Item(24, 13).Value = 31.06110666666667
$ws.Cells.Item(24, 14).Value = 93.18332000000001
$ws.Cells.Item(24, 15).Value = 0.4341659082002267
$ws.Cells.Item(24, 16).Value = 0.4341659082002268
$ws.Cells.Item(24, 17).Value = 183.1676152095912
$ws.Cells.Item(24, 18).Value = 1648.50853688632
$ws.Cells.Item(24, 19).Value = 0.05256481864297925
$ws.Cells.Item(24, 20).Value = 0.05256481864297926

$ws.Cells.Item(25, 1).Value = "M2"
$ws.Cells.Item(25, 2).Value = "Agrn"
$ws.Cells.Item(25, 3).Value = "Atp1a3"
$ws.Cells.Item(25, 4).Value = "sCs"
$ws.Cells.Item(25, 5).Value = 3
$ws.Cells.Item(25, 6).Value = 1
$ws.Cells.Item(25, 7).Value = 5.897008666666667
$ws.Cells.Item(25, 8).Value = 17.691026
$ws.Cells.Item(25, 9).Value = 0.1210708110659339
$ws.Cells.Item(25, 10).Value = 0.1210708110659339
$ws.Cells.Item(25, 11).Value = 3
$ws.Cells.Item(25, 12).Value = 1
$ws.Cells.Item(25, 13).Value = 0.07433466666666666
$ws.Cells.Item(25, 14).Value = 0.223004
$ws.Cells.Item(25, 15).Value = 0.001039035035372032
$ws.Cells.Item(25, 16).Value = 0.001039035035372032
$ws.Cells.Item(25, 17).Value = 0.4383521735671111
$ws.Cells.Item(25, 18).Value = 3.945169562104
$ws.Cells.Item(25, 19).Value = 0.0001257968144584132
$ws.Cells.Item(25, 20).Value = 0.0001257968144584132

$ws.Cells.Item(26, 1).Value = "Neutro"
$ws.Cells.Item(26, 2).Value = "Agrn"
$ws.Cells.Item(26, 3).Value = "Atp1a3"
$ws.Cells.Item(26, 4).Value = "ECs"
$ws.Cells.Item(26, 5).Value = 3
$ws.Cells.Item(26, 6).Value = 1
$ws.Cells.Item(26, 7).Value = 5.159083666666667
$ws.Cells.Item(26, 8).Value = 15.477251
$ws.Cells.Item(26, 9).Value = 0.1059205572159035
$ws.Cells.Item(26, 10).Value = 0.1059205572159035
$ws.Cells.Item(26, 11).Value = 2
$ws.Cells.Item(26, 12).Value = 0.6666666666666666
$ws.Cells.Item(26, 13).Value = 0.012411
$ws.Cells.Item(26, 14).Value = 0.037233
$ws.Cells.Item(26, 15).Value = 0.0001734784643863198
$ws.Cells.Item(26, 16).Value = 0.0001734784643863198
$ws.Cells.Item(26, 17).Value = 0.06402938738699999
$ws.Cells.Item(26, 18).Value = 0.576264486483
$ws.Cells.Item(26, 19).Value = 0.00001837493561275826
$ws.Cells.Item(26, 20).Value = 0.00001837493561275827

$ws.Cells.Item(27, 1).Value = "Neutro"
$ws.Cells.Item(27, 2).Value = "Agrn"
$ws.Cells.Item(27, 3).Value = "Atp1a3"
$ws.Cells.Item(27, 4).Value = "FAPs"
$ws.Cells.Item(27, 5).Value = 3
$ws.Cells.Item(27, 6).Value = 1
$ws.Cells.Item(27, 7).Value = 5.159083666666667
$ws.Cells.Item(27, 8).Value = 15.477251
$ws.Cells.Item(27, 9).Value = 0.1059205572159035
$ws.Cells.Item(27, 10).Value = 0.1059205572159035
$ws.Cells.Item(27, 11).Value = 3
$ws.Cells.Item(27, 12).Value = 1
$ws.Cells.Item(27, 13).Value = 0.08803833333333333
$ws.Cells.Item(27, 14).Value = 0.264115
$ws.Cells.Item(27, 15).Value = 0.001230582134702894
$ws.Cells.Item(27, 16).Value = 0.001230582134702894
$ws.Cells.Item(27, 17).Value = 0.4541971275405555
$ws.Cells.Item(27, 18).Value = 4.087774147865
$ws.Cells.Item(27, 19).Value = 0.0001303439454076665
$ws.Cells.Item(27, 20).Value = 0.0001303439454076666

$ws.Cells.Item(28, 1).Value = "Neutro"
$ws.Cells.Item(28, 2).Value = "Agrn"
$ws.Cells.Item(28, 3).Value = "Atp1a3"
$ws.Cells.Item(28, 4).Value = "M1"
$ws.Cells.Item(28, 5).Value = 3
$ws.Cells.Item(28, 6).Value = 1
$ws.Cells.Item(28, 7).Value = 5.159083666666667
$ws.Cells.Item(28, 8).Value = 15.477251
$ws.Cells.Item(28, 9).Value = 0.1059205572159035
$ws.Cells.Item(28, 10).Value = 0.1059205572159035
$ws.Cells.Item(28, 11).Value = 3
$ws.Cells.Item(28, 12).Value = 1
$ws.Cells.Item(28, 13).Value = 20.739774
$ws.Cells.Item(28, 14).Value = 62.219322
$ws.Cells.Item(28, 15).Value = 0.2898963939440272
$ws.Cells.Item(28, 16).Value = 0.2898963939440272
$ws.Cells.Item(28, 17).Value = 106.998229293758
$ws.Cells.Item(28, 18).Value = 962.9840636438219
$ws.Cells.Item(28, 19).Value = 0.03070598758143242
$ws.Cells.Item(28, 20).Value = 0.03070598758143244

$ws.Cells.Item(29, 1).Value = "Neutro"
$ws.Cells.Item(29, 2).Value = "Agrn"
$ws.Cells.Item(29, 3).Value = "Atp1a3"
$ws.Cells.Item(29, 4).Value = "M2"
$ws.Cells.Item(29, 5).Value = 3
$ws.Cells.Item(29, 6).Value = 1
$ws.Cells.Item(29, 7).Value = 5.159083666666667
$ws.Cells.Item(29, 8).Value = 15.477251
$ws.Cells.Item(29, 9).Value = 0.1059205572159035
$ws.Cells.Item(29, 10).Value = 0.1059205572159035
$ws.Cells.Item(29, 11).Value = 3
$ws.Cells.Item(29, 12).Value = 1
$ws.Cells.Item(29, 13).Value = 19.56635666666667
$ws.Cells.Item(29, 14).Value = 58.69907
$ws.Cells.Item(29, 15).Value = 0.2734946022212847
$ws.Cells.Item(29, 16).Value = 0.2734946022212847
$ws.Cells.Item(29, 17).Value = 100.9444710951744
$ws.Cells.Item(29, 18).Value = 908.5002398565699
$ws.Cells.Item(29, 19).Value = 0.02896870066282035
$ws.Cells.Item(29, 20).Value = 0.02896870066282035

$ws.Cells.Item(30, 1).Value = "Neutro"
$ws.Cells.Item(30, 2).Value = "Agrn"
$ws.Cells.Item(30, 3).Value = "Atp1a3"
$ws.Cells.Item(30, 4).Value = "Neutro"
$ws.Cells.Item(30, 5).Value = 3
$ws.Cells.Item(30, 6).Value = 1
$ws.Cells.Item(30, 7).Value = 5.159083666666667
$ws.Cells.Item(30, 8).Value = 15.477251
$ws.Cells.Item(30, 9).Value = 0.1059205572159035
$ws.Cells.Item(30, 10).Value = 0.1059205572159035
$ws.Cells.Item(30, 11).Value = 3
$ws.Cells.Item(30, 12).Value = 1
$ws.Cells.Item(30, 13).Value = 31.06110666666667
$ws.Cells.Item(30, 14).Value = 93.18332000000001
$ws.Cells.Item(30, 15).Value = 0.4341659082002267
$ws.Cells.Item(30, 16).Value = 0.4341659082002268
$ws.Cells.Item(30, 17).Value = 160.2468480725911
$ws.Cells.Item(30, 18).Value = 1442.22163265332
$ws.Cells.Item(30, 19).Value = 0.04598709492071681
$ws.Cells.Item(30, 20).Value = 0.04598709492071684

$ws.Cells.Item(31, 1).Value = "Neutro"
$ws.Cells.Item(31, 2).Value = "Agrn"
$ws.Cells.Item(31, 3).Value = "Atp1a3"
$ws.Cells.Item(31, 4).Value = "sCs"
$ws.Cells.Item(31, 5).Value = 3
$ws.Cells.Item(31, 6).Value = 1
$ws.Cells.Item(31, 7).Value = 5.159083666666667
$ws.Cells.Item(31, 8).Value = 15.477251
$ws.Cells.Item(31, 9).Value = 0.1059205572159035
$ws.Cells.Item(31, 10).Value = 0.1059205572159035
$ws.Cells.Item(31, 11).Value = 3
$ws.Cells.Item(31, 12).Value = 1
$ws.Cells.Item(31, 13).Value = 0.07433466666666666
$ws.Cells.Item(31, 14).Value = 0.223004
$ws.Cells.Item(31, 15).Value = 0.001039035035372032
$ws.Cells.Item(31, 16).Value = 0.001039035035372032
$ws.Cells.Item(31, 17).Value = 0.3834987646671111
$ws.Cells.Item(31, 18).Value = 3.451488882004
$ws.Cells.Item(31, 19).Value = 0.0001100551699134516
$ws.Cells.Item(31, 20).Value = 0.0001100551699134516

$ws.Cells.Item(32, 1).Value = "sCs"
$ws.Cells.Item(32, 2).Value = "Agrn"
$ws.Cells.Item(32, 3).Value = "Atp1a3"
$ws.Cells.Item(32, 4).Value = "ECs"
$ws.Cells.Item(32, 5).Value = 3
$ws.Cells.Item(32, 6).Value = 1
$ws.Cells.Item(32, 7).Value = 12.604947
$ws.Cells.Item(32, 8).Value = 37.814841
$ws.Cells.Item(32, 9).Value = 0.2587907264507627
$ws.Cells.Item(32, 10).Value = 0.2587907264507627
$ws.Cells.Item(32, 11).Value = 2
$ws.Cells.Item(32, 12).Value = 0.6666666666666666
$ws.Cells.Item(32, 13).Value = 0.012411
$ws.Cells.Item(32, 14).Value = 0.037233
$ws.Cells.Item(32, 15).Value = 0.0001734784643863198
$ws.Cells.Item(32, 16).Value = 0.0001734784643863198
$ws.Cells.Item(32, 17).Value = 0.156439997217
$ws.Cells.Item(32, 18).Value = 1.407959974953
$ws.Cells.Item(32, 19).Value = 0.00004489461782209846
$ws.Cells.Item(32, 20).Value = 0.00004489461782209848

$ws.Cells.Item(33, 1).Value = "sCs"
$ws.Cells.Item(33, 2).Value = "Agrn"
$ws.Cells.Item(33, 3).Value = "Atp1a3"
$ws.Cells.Item(33, 4).Value = "FAPs"
$ws.Cells.Item(33, 5).Value = 3
$ws.Cells.Item(33, 6).Value = 1
$ws.Cells.Item(33, 7).Value = 12.604947
$ws.Cells.Item(33, 8).Value = 37.814841
$ws.Cells.Item(33, 9).Value = 0.2587907264507627
$ws.Cells.Item(33, 10).Value = 0.2587907264507627
$ws.Cells.Item(33, 11).Value = 3
$ws.Cells.Item(33, 12).Value = 1
$ws.Cells.Item(33, 13).Value = 0.08803833333333333
$ws.Cells.Item(33, 14).Value = 0.264115
$ws.Cells.Item(33, 15).Value = 0.001230582134702894
$ws.Cells.Item(33, 16).Value = 0.001230582134702894
$ws.Cells.Item(33, 17).Value = 1.109718525635
$ws.Cells.Item(33, 18).Value = 9.987466730714999
$ws.Cells.Item(33, 19).Value = 0.0003184632445970922
$ws.Cells.Item(33, 20).Value = 0.0003184632445970923

$ws.Cells.Item(34, 1).Value = "sCs"
$ws.Cells.Item(34, 2).Value = "Agrn"
$ws.Cells.Item(34, 3).Value = "Atp1a3"
$ws.Cells.Item(34, 4).Value = "M1"
$ws.Cells.Item(34, 5).Value = 3
$ws.Cells.Item(34, 6).Value = 1
$ws.Cells.Item(34, 7).Value = 12.604947
$ws.Cells.Item(34, 8).Value = 37.814841
$ws.Cells.Item(34, 9).Value = 0.2587907264507627
$ws.Cells.Item(34, 10).Value = 0.2587907264507627
$ws.Cells.Item(34, 11).Value = 3
$ws.Cells.Item(34, 12).Value = 1
$ws.Cells.Item(34, 13).Value = 20.739774
$ws.Cells.Item(34, 14).Value = 62.219322
$ws.Cells.Item(34, 15).Value = 0.2898963939440272
$ws.Cells.Item(34, 16).Value = 0.2898963939440272
$ws.Cells.Item(34, 17).Value = 261.423752061978
$ws.Cells.Item(34, 18).Value = 2352.813768557802
$ws.Cells.Item(34, 19).Value = 0.07502249838423126
$ws.Cells.Item(34, 20).Value = 0.07502249838423129

$ws.Cells.Item(35, 1).Value = "sCs"
$ws.Cells.Item(35, 2).Value = "Agrn"
$ws.Cells.Item(35, 3).Value = "Atp1a3"
$ws.Cells.Item(35, 4).Value = "M2"
$ws.Cells.Item(35, 5).Value = 3
$ws.Cells.Item(35, 6).Value = 1
$ws.Cells.Item(35, 7).Value = 12.604947
$ws.Cells.Item(35, 8).Value = 37.814841
$ws.Cells.Item(35, 9).Value = 0.2587907264507627
$ws.Cells.Item(35, 10).Value = 0.2587907264507627
$ws.Cells.Item(35, 11).Value = 3
$ws.Cells.Item(35, 12).Value = 1
$ws.Cells.Item(35, 13).Value = 19.56635666666667
$ws.Cells.Item(35, 14).Value = 58.69907
$ws.Cells.Item(35, 15).Value = 0.2734946022212847
$ws.Cells.Item(35, 16).Value = 0.2734946022212847
$ws.Cells.Item(35, 17).Value = 246.63288876643
$ws.Cells.Item(35, 18).Value = 2219.69599889787
$ws.Cells.Item(35, 19).Value = 0.07077786678920864
$ws.Cells.Item(35, 20).Value = 0.07077786678920865

$ws.Cells.Item(36, 1).Value = "sCs"
$ws.Cells.Item(36, 2).Value = "Agrn"
$ws.Cells.Item(36, 3).Value = "Atp1a3"
$ws.Cells.Item(36, 4).Value = "Neutro"
$ws.Cells.Item(36, 5).Value = 3
$ws.Cells.Item(36, 6).Value = 1
$ws.Cells.Item(36, 7).Value = 12.604947
$ws.Cells.Item(36, 8).Value = 37.814841
$ws.Cells.Item(36, 9).Value = 0.2587907264507627
$ws.Cells.Item(36, 10).Value = 0.2587907264507627
$ws.Cells.Item(36, 11).Value = 3
$ws.Cells.Item(36, 12).Value = 1
$ws.Cells.Item(36, 13).Value = 31.06110666666667
$ws.Cells.Item(36, 14).Value = 93.18332000000001
$ws.Cells.Item(36, 15).Value = 0.4341659082002267
$ws.Cells.Item(36, 16).Value = 0.4341659082002268
$ws.Cells.Item(36, 17).Value = 391.5236032946801
$ws.Cells.Item(36, 18).Value = 3523.712429652121
$ws.Cells.Item(36, 19).Value = 0.1123581107832918
$ws.Cells.Item(36, 20).Value = 0.1123581107832919

$ws.Cells.Item(37, 1).Value = "sCs"
$ws.Cells.Item(37, 2).Value = "Agrn"
$ws.Cells.Item(37, 3).Value = "Atp1a3"
$ws.Cells.Item(37, 4).Value = "sCs"
$ws.Cells.Item(37, 5).Value = 3
$ws.Cells.Item(37, 6).Value = 1
$ws.Cells.Item(37, 7).Value = 12.604947
$ws.Cells.Item(37, 8).Value = 37.814841
$ws.Cells.Item(37, 9).Value = 0.2587907264507627
$ws.Cells.Item(37, 10).Value = 0.2587907264507627
$ws.Cells.Item(37, 11).Value = 3
$ws.Cells.Item(37, 12).Value = 1
$ws.Cells.Item(37, 13).Value = 0.07433466666666666
$ws.Cells.Item(37, 14).Value = 0.223004
$ws.Cells.Item(37, 15).Value = 0.001039035035372032
$ws.Cells.Item(37, 16).Value = 0.001039035035372032
$ws.Cells.Item(37, 17).Value = 0.936984533596
$ws.Cells.Item(37, 18).Value = 8.432860802363999
$ws.Cells.Item(37, 19).Value = 0.000268892631611722
$ws.Cells.Item(37, 20).Value = 0.0002688926316117221

Write-Host "Applied Agrn-Atp1a3 NATMI re-run update (Dr Hou advice)"
